$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing phone numbers (E2, E3, E4, E10)
$ws.Range("E2").Value = 523301800
$ws.Range("E3").Value = 523301802
$ws.Range("E4").Value = 523301801
$ws.Range("E10").Value = 528827064

# Add new row 11 - duplicate of row 2 formatting/content pattern, with new name & phone
$ws.Range("A11").Value = "ינון הדר12"
$ws.Range("B11").Value = "בני דוד עלי"
$ws.Range("C11").Value = "מלווה"
$ws.Range("D11").Value = "aviadvcux@Gmail.com"
$ws.Range("E11").Value = 523301803
$ws.Range("F11").Value = "מרכז"

# Copy formatting from row 2 (A2:U2) to row 11 (A11:U11)
$ws.Range("A2:U2").Copy()
$ws.Range("A11:U11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-set values after paste (paste special formats only, but ensure values are correct)
$ws.Range("A11").Value = "ינון הדר12"
$ws.Range("B11").Value = "בני דוד עלי"
$ws.Range("C11").Value = "מלווה"
$ws.Range("D11").Value = "aviadvcux@Gmail.com"
$ws.Range("E11").Value = 523301803
$ws.Range("F11").Value = "מרכז"

$ws.Range("E11").Select()

$wb.Save()
